# Re-run SGNN to annotate dialog acts following clean up work to the
# original transcripts. This updates the DAMSLTag (column I) and
# DialogAct (column J) values for a set of rows in Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @(
    @{ Row = 2;   I = "b";  J = "Acknowledge (Backchannel)" },
    @{ Row = 39;  I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 49;  I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 50;  I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 59;  I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 61;  I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 75;  I = "sv"; J = "Statement-opinion" },
    @{ Row = 77;  I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 79;  I = "%";  J = "Uninterpretable" },
    @{ Row = 91;  I = "ba"; J = "Appreciation" },
    @{ Row = 98;  I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 102; I = "sv"; J = "Statement-opinion" },
    @{ Row = 132; I = "sv"; J = "Statement-opinion" },
    @{ Row = 144; I = "sv"; J = "Statement-opinion" },
    @{ Row = 152; I = "sv"; J = "Statement-opinion" },
    @{ Row = 169; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 179; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 188; I = "sv"; J = "Statement-opinion" },
    @{ Row = 238; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 256; I = "ba"; J = "Appreciation" },
    @{ Row = 263; I = "sv"; J = "Statement-opinion" },
    @{ Row = 271; I = "sv"; J = "Statement-opinion" },
    @{ Row = 275; I = "sv"; J = "Statement-opinion" },
    @{ Row = 276; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 279; I = "aa"; J = "Agree/Accept" },
    @{ Row = 282; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 289; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 308; I = "sv"; J = "Statement-opinion" }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.I
    $ws.Cells.Item($u.Row, 10).Value = $u.J
}
